$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.854221333333333
$ws.Range("H2").Value = 5.562664
$ws.Range("I2").Value = 0.03416002559055492
$ws.Range("J2").Value = 0.03416002559055492
$ws.Range("M2").Value = 211.980367
$ws.Range("N2").Value = 635.9411009999999
$ws.Range("O2").Value = 0.9885149156420702
$ws.Range("P2").Value = 0.9885149156420702
$ws.Range("Q2").Value = 393.0585187392293
$ws.Range("R2").Value = 3537.526668653063
$ws.Range("S2").Value = 0.03376769481497836
$ws.Range("T2").Value = 0.03376769481497836

$ws.Range("G3").Value = 1.854221333333333
$ws.Range("H3").Value = 5.562664
$ws.Range("I3").Value = 0.03416002559055492
$ws.Range("J3").Value = 0.03416002559055492
$ws.Range("O3").Value = 0.003992992409159323
$ws.Range("P3").Value = 0.003992992409159324
$ws.Range("Q3").Value = 1.587714719167111
$ws.Range("R3").Value = 14.289432472504
$ws.Range("S3").Value = 0.000136400722879774
$ws.Range("T3").Value = 0.000136400722879774

$ws.Range("G4").Value = 1.854221333333333
$ws.Range("H4").Value = 5.562664
$ws.Range("I4").Value = 0.03416002559055492
$ws.Range("J4").Value = 0.03416002559055492
$ws.Range("O4").Value = 0.007492091948770576
$ws.Range("P4").Value = 0.007492091948770576
$ws.Range("Q4").Value = 2.979045148478222
$ws.Range("R4").Value = 26.811406336304
$ws.Range("S4").Value = 0.0002559300526967934
$ws.Range("T4").Value = 0.0002559300526967934

$ws.Range("I5").Value = 0.8311547934421808
$ws.Range("J5").Value = 0.8311547934421808
$ws.Range("M5").Value = 211.980367
$ws.Range("N5").Value = 635.9411009999999
$ws.Range("O5").Value = 0.9885149156420702
$ws.Range("P5").Value = 0.9885149156420702
$ws.Range("Q5").Value = 9563.589789690395
$ws.Range("R5").Value = 86072.30810721355
$ws.Range("S5").Value = 0.8216089105249996
$ws.Range("T5").Value = 0.8216089105249996

$ws.Range("I6").Value = 0.8311547934421808
$ws.Range("J6").Value = 0.8311547934421808
$ws.Range("O6").Value = 0.003992992409159323
$ws.Range("P6").Value = 0.003992992409159324
$ws.Range("S6").Value = 0.003318794781051013
$ws.Range("T6").Value = 0.003318794781051013

$ws.Range("I7").Value = 0.8311547934421808
$ws.Range("J7").Value = 0.8311547934421808
$ws.Range("O7").Value = 0.007492091948770576
$ws.Range("P7").Value = 0.007492091948770576
$ws.Range("S7").Value = 0.006227088136130234
$ws.Range("T7").Value = 0.006227088136130234

$ws.Range("I8").Value = 0.1346851809672642
$ws.Range("J8").Value = 0.1346851809672642
$ws.Range("M8").Value = 211.980367
$ws.Range("N8").Value = 635.9411009999999
$ws.Range("O8").Value = 0.9885149156420702
$ws.Range("P8").Value = 0.9885149156420702
$ws.Range("Q8").Value = 1549.739990293079
$ws.Range("R8").Value = 13947.65991263771
$ws.Range("S8").Value = 0.1331383103020922
$ws.Range("T8").Value = 0.1331383103020922

$ws.Range("I9").Value = 0.1346851809672642
$ws.Range("J9").Value = 0.1346851809672642
$ws.Range("O9").Value = 0.003992992409159323
$ws.Range("P9").Value = 0.003992992409159324
$ws.Range("S9").Value = 0.0005377969052285358
$ws.Range("T9").Value = 0.0005377969052285359

$ws.Range("I10").Value = 0.1346851809672642
$ws.Range("J10").Value = 0.1346851809672642
$ws.Range("O10").Value = 0.007492091948770576
$ws.Range("P10").Value = 0.007492091948770576
$ws.Range("S10").Value = 0.001009073759943548
$ws.Range("T10").Value = 0.001009073759943548
